$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two data rows (old rows 2 and 3) are removed; remaining rows shift up.
$ws.Rows("2:3").Delete()

# Append 12 new data rows (new rows 20-31) with freshly captured gyroscope samples.
$ws.Range("A20").Value = -1.937935614804442
$ws.Range("B20").Value = -0.1515454994429453
$ws.Range("C20").Value = -1.15992745128265
$ws.Range("A21").Value = 2.910339694743533
$ws.Range("B21").Value = -5.836961925576653
$ws.Range("C21").Value = -4.45155078828873
$ws.Range("A22").Value = 3.819924138400946
$ws.Range("B22").Value = -3.973247524795179
$ws.Range("C22").Value = -0.4456482636272369
$ws.Range("A23").Value = 1.989335611325833
$ws.Range("B23").Value = -0.2693817703002244
$ws.Range("C23").Value = 0.09849660987155473
$ws.Range("A24").Value = -3.755865255627084
$ws.Range("B24").Value = 0.09460148078591657
$ws.Range("C24").Value = 3.946894074251968
$ws.Range("A25").Value = -2.889570736010117
$ws.Range("B25").Value = 9.743136939652585
$ws.Range("C25").Value = 0.5475007984615461
$ws.Range("A26").Value = -1.842574656009669
$ws.Range("B26").Value = 4.510396480560456
$ws.Range("C26").Value = -1.828126966953318
$ws.Range("A27").Value = 2.88307266279109
$ws.Range("B27").Value = -8.930636825911506
$ws.Range("C27").Value = -1.004032475139075
$ws.Range("A28").Value = 4.052329929460007
$ws.Range("B28").Value = -3.608018367662238
$ws.Range("C28").Value = -3.160855819325922
$ws.Range("A29").Value = -1.29231422538059
$ws.Range("B29").Value = 1.658297274637648
$ws.Range("C29").Value = -0.279984981641821
$ws.Range("A30").Value = -4.269718651377847
$ws.Range("B30").Value = -2.122582794056059
$ws.Range("C30").Value = -2.14775515915062
$ws.Range("A31").Value = -2.918748300009885
$ws.Range("B31").Value = 1.177203515254081
$ws.Range("C31").Value = -6.725167151984817
